# Updated Test cases and bug metrics
# All test cases for iteration PASSED -> mark every bug row as "Solved"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bug Metrics")

# Column I (Status) for rows 8-17 flips from "Unsolved" to "Solved"
$ws.Range("I8:I17").Value = "Solved"

# Restore the view to show the top of the table with J9 selected
$ws.Activate() | Out-Null
$ws.Range("J9").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 2
